$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vlan_params")

# Add a new row for the is_dhcp check, mirroring the existing rows.
# Seed A10 by copying the format/type from A9, then set its text.
$ws.Range("A9").Copy($ws.Range("A10"))
$ws.Range("A10").Value = "is_dhcp"

# B3 already holds the text value "True" (not a real Boolean). Copy it
# straight into B10 so the new cell keeps the same plain-text type and
# style instead of Excel auto-converting a typed "True" into a Boolean.
$ws.Range("B3").Copy($ws.Range("B10"))

# Update selection to reflect where the cursor ended up after the edit
$ws.Range("B11").Select()

$wb.Save()
